# Re-sync "Students", "Additional", "Points" and "Score" sheets after
# grading corrections for Student A / Student B / Student C
# (fixes an autosave/save/load round-trip bug).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Students: Task 2 / Task 4 points were corrected for three students,
# which shifts their Grade and Total columns.
# ---------------------------------------------------------------------
$students = $wb.Worksheets.Item("Students")

# Student A (row 2): Task2 8 -> 10, Task4 6 -> 3.5
$students.Range("F2").Value = 10
$students.Range("H2").Value = 3.5
$students.Range("B2").Value = "2-"
$students.Range("D2").Value = "22.5 (72.6%)"

# Student B (row 3): Task4 5 -> 3
$students.Range("H3").Value = 3
$students.Range("B3").Value = "2+"
$students.Range("D3").Value = "25.0 (80.6%)"

# Student C (row 4): Task2 1 -> 8
$students.Range("F4").Value = 8
$students.Range("B4").Value = "2+"
$students.Range("D4").Value = "25.0 (80.6%)"

# ---------------------------------------------------------------------
# Additional: the "passed" flag flips now that the corrected totals no
# longer clear the bar.
# ---------------------------------------------------------------------
$additional = $wb.Worksheets.Item("Additional")
$additional.Range("A2").Value = $false

# ---------------------------------------------------------------------
# Points: per-total breakdown table, resorted/regrouped for the new
# totals (22.5 for Student A, 25.0 shared by Student B & Student C).
# ---------------------------------------------------------------------
$points = $wb.Worksheets.Item("Points")

$points.Range("A5").Value = 19
$points.Range("B5").NumberFormat = "@"
$points.Range("B5").Value = "3"
$points.Range("D5").Value = "Student H"

$points.Range("A6").Value = 20
$points.Range("B6").Value = "3+"
$points.Range("D6").Value = "Student I"

$points.Range("A7").Value = 22
$points.Range("B7").Value = "2-"
$points.Range("D7").Value = "Student J"

$points.Range("A8").Value = 22.5
$points.Range("D8").Value = "Student A"

$points.Range("A9").Value = 25
$points.Range("B9").Value = "2+"
$points.Range("C9").Value = 2
$points.Range("D9").Value = "Student B, Student C"

$points.Range("C10").Value = 1
$points.Range("D10").Value = "Student F"

# ---------------------------------------------------------------------
# Score: previously grouped by the bare grade number (a bug); now
# grouped by the precise grade (with +/- modifiers), so the table grows
# from 6 data rows to 8.
# ---------------------------------------------------------------------
$score = $wb.Worksheets.Item("Score")

$score.Range("A2").Value = "1-"
$score.Range("B2").Value = 1
$score.Range("C2").Value = "Student F"

$score.Range("A3").Value = "2+"
$score.Range("B3").Value = 2
$score.Range("C3").Value = "Student B, Student C"

$score.Range("A4").Value = "2-"
$score.Range("B4").Value = 2
$score.Range("C4").Value = "Student A, Student J"

$score.Range("A5").NumberFormat = "@"
$score.Range("A5").Value = "3"
$score.Range("B5").Value = 1
$score.Range("C5").Value = "Student H"

$score.Range("A6").Value = "3+"
$score.Range("B6").Value = 1
$score.Range("C6").Value = "Student I"

$score.Range("A7").Value = "4-"
$score.Range("B7").Value = 1
$score.Range("C7").Value = "Student G"

# A8/A9 are brand-new cells below the sheet's old A1:C7 extent, so they
# don't inherit the bordered/bold "grade label" look of A2:A7 for free;
# copy that formatting across before writing the new values.
$score.Range("A2").Copy()
$score.Range("A8").PasteSpecial(-4122)
$score.Range("A8").Value = "5-"
$score.Range("B8").Value = 1
$score.Range("C8").Value = "Student D"

$score.Range("A2").Copy()
$score.Range("A9").PasteSpecial(-4122)
$score.Range("A9").NumberFormat = "@"
$score.Range("A9").Value = "6"
$score.Range("B9").Value = 1
$score.Range("C9").Value = "Student E"
